# Projektvorbereitung / Zeitplan u. Präsi update
# 1) Add a bold "Milestones" heading paragraph to the
#    "Projekt- und Zeitplan (Milestones)" slide's bullet list.
# 2) Append two new "Szenario" slides (Title + Content layout) at the end
#    of the deck.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 5 ("Projekt- und Zeitplan (Milestones)") - insert a bold, un-bulleted
#    "Milestones" paragraph above the existing numbered milestone list.
# ---------------------------------------------------------------------------
$milestoneSlide = $p.Slides.Item(5)
$contentShape = $milestoneSlide.Shapes.Item(2)
$contentTf = $contentShape.TextFrame

# Insert a new first paragraph carrying the text, then strip its bullet /
# indent so it reads as a plain bold heading line above the numbered list.
$contentTf.TextRange.InsertBefore("Milestones`r")
$headingPara = $contentTf.TextRange.Paragraphs(1, 1)

$rulerLvl1 = $contentTf.Ruler.Levels.Item(1)
$rulerLvl1.LeftMargin = 0
$rulerLvl1.FirstMargin = 0

$headingPara.ParagraphFormat.Bullet.Visible = 0
$headingPara.Font.Bold = 1

# ---------------------------------------------------------------------------
# 2) Append two new "Szenario" slides. Duplicating an existing Title+Content
#    slide (and moving the copy to the end) keeps the placeholder naming,
#    ids, language and creation-id metadata consistent with the rest of the
#    deck; we then overwrite the title and clear the body placeholder.
# ---------------------------------------------------------------------------
$templateSlide = $p.Slides.Item(8)

for ($i = 0; $i -lt 2; $i++) {
    $newSlide = $templateSlide.Duplicate()
    $newSlide.MoveTo($p.Slides.Count)

    $newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Szenario"

    $bodyTf = $newSlide.Shapes.Item(2).TextFrame
    # Remove every paragraph but the last (blank) one so the placeholder
    # ends up empty, exactly like a freshly-inserted content placeholder.
    for ($j = 0; $j -lt 5; $j++) {
        $bodyTf.TextRange.Paragraphs(1, 1).Delete()
    }
}
